$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scheme names ("Gaussian-Quadrature" and three spiral schemes) were
# inserted into the ordered list of sampling schemes, right after
# "Ring Perpendicular to TD" and before "NoRotation-tilt60deg".
# Existing rows 10-16 (A index 8-14) shift down to make room; their
# labels now correspond to the schemes that used to sit after the
# inserted block, continuing in the same overall ordered sequence.
$shiftedLabels = @(
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt"
)

for ($i = 0; $i -lt $shiftedLabels.Length; $i++) {
    $row = 10 + $i
    $ws.Cells.Item($row, 2).Value = $shiftedLabels[$i]
}

# Three new trailing rows for the schemes pushed off the end of the
# original table (the HexGrid variants).
$trailingLabels = @(
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)

for ($i = 0; $i -lt $trailingLabels.Length; $i++) {
    $row = 17 + $i
    $a = 15 + $i
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $trailingLabels[$i]
    for ($col = 3; $col -le 16; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
    # Match the formatting already used on the rest of column A (bold,
    # bordered, centered/top-aligned header-style cell).
    $ws.Cells.Item(16, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
}
